$d = $word.ActiveDocument

$d.Content.Find.Execute("40×29=", $true, $false, $false, $false, $false, $true, 1, $false, "51×84=", 2) | Out-Null
$d.Content.Find.Execute("25×22=", $true, $false, $false, $false, $false, $true, 1, $false, "76×59=", 2) | Out-Null
$d.Content.Find.Execute("74×59=", $true, $false, $false, $false, $false, $true, 1, $false, "40×62=", 2) | Out-Null
$d.Content.Find.Execute("72×23=", $true, $false, $false, $false, $false, $true, 1, $false, "35×65=", 2) | Out-Null
$d.Content.Find.Execute("11×20=", $true, $false, $false, $false, $false, $true, 1, $false, "14×65=", 2) | Out-Null
$d.Content.Find.Execute("59×12=", $true, $false, $false, $false, $false, $true, 1, $false, "61×23=", 2) | Out-Null
$d.Content.Find.Execute("23×72=", $true, $false, $false, $false, $false, $true, 1, $false, "34×79=", 2) | Out-Null
$d.Content.Find.Execute("18×60=", $true, $false, $false, $false, $false, $true, 1, $false, "40×94=", 2) | Out-Null
$d.Content.Find.Execute("57×60=", $true, $false, $false, $false, $false, $true, 1, $false, "73×76=", 2) | Out-Null
$d.Content.Find.Execute("22×80=", $true, $false, $false, $false, $false, $true, 1, $false, "27×47=", 2) | Out-Null
$d.Content.Find.Execute("53×73=", $true, $false, $false, $false, $false, $true, 1, $false, "76×99=", 2) | Out-Null
$d.Content.Find.Execute("41×82=", $true, $false, $false, $false, $false, $true, 1, $false, "22×50=", 2) | Out-Null
$d.Content.Find.Execute("24×29=", $true, $false, $false, $false, $false, $true, 1, $false, "89×74=", 2) | Out-Null
$d.Content.Find.Execute("14×90=", $true, $false, $false, $false, $false, $true, 1, $false, "73×77=", 2) | Out-Null
$d.Content.Find.Execute("60×83=", $true, $false, $false, $false, $false, $true, 1, $false, "36×24=", 2) | Out-Null
$d.Content.Find.Execute("75×99=", $true, $false, $false, $false, $false, $true, 1, $false, "49×92=", 2) | Out-Null
$d.Content.Find.Execute("24×96=", $true, $false, $false, $false, $false, $true, 1, $false, "81×37=", 2) | Out-Null
$d.Content.Find.Execute("67×61=", $true, $false, $false, $false, $false, $true, 1, $false, "24×14=", 2) | Out-Null
$d.Content.Find.Execute("25×80=", $true, $false, $false, $false, $false, $true, 1, $false, "46×74=", 2) | Out-Null
$d.Content.Find.Execute("69×53=", $true, $false, $false, $false, $false, $true, 1, $false, "48×44=", 2) | Out-Null
$d.Content.Find.Execute("15×13=", $true, $false, $false, $false, $false, $true, 1, $false, "81×42=", 2) | Out-Null
$d.Content.Find.Execute("94×24=", $true, $false, $false, $false, $false, $true, 1, $false, "24×72=", 2) | Out-Null
$d.Content.Find.Execute("32×48=", $true, $false, $false, $false, $false, $true, 1, $false, "18×14=", 2) | Out-Null
$d.Content.Find.Execute("51×60=", $true, $false, $false, $false, $false, $true, 1, $false, "15×35=", 2) | Out-Null
$d.Content.Find.Execute("78×21=", $true, $false, $false, $false, $false, $true, 1, $false, "32×91=", 2) | Out-Null
